$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook's last data row is row 77 (a single intraday print for
# 2024-06-12 14:53...). The R script re-ran and now reports a clean
# daily bar for 2024-06-11 (new row 77) followed by the previously
# recorded 2024-06-12 bar, but re-stamped to a clean 07:00 timestamp
# (pushed down to row 78). Insert a fresh row so the old row 77 slides
# down to row 78, then (re)write both rows' values.
$ws.Rows.Item(77).Insert()

# New row 77: 2024-06-11 data
$ws.Range("A77").Value = 45454.2916666667
$ws.Range("B77").Value = 4500
$ws.Range("C77").Value = 3
$ws.Range("D77").Value = 2.99000000953674
$ws.Range("E77").Value = 2.99000000953674
$ws.Range("F77").Value = 3
$ws.Range("G77").NumberFormat = "@"
$ws.Range("G77").Value = "3"
$ws.Range("G77").Style = "Normal"
$ws.Range("H77").Value = "ESPE.MI"

# Row 78: same values the old row 77 had, just a cleaned-up date stamp
$ws.Range("A78").Value = 45455.2916666667
$ws.Range("B78").Value = 18000
$ws.Range("C78").Value = 3.27999997138977
$ws.Range("D78").Value = 2.99000000953674
$ws.Range("E78").Value = 2.99000000953674
$ws.Range("F78").Value = 2.99000000953674
$ws.Range("G78").NumberFormat = "@"
$ws.Range("G78").Value = "2.99000000953674"
$ws.Range("G78").Style = "Normal"
$ws.Range("H78").Value = "ESPE.MI"
